# Rename the "Speaker" sheet to "Name" and make it the active/selected sheet
# (previously "Transition" was the active sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Speaker")
$ws.Name = "Name"
$ws.Activate()
